$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final table (rows 2-7) for columns A..T, following the updated
# Sending-cluster x Target-cluster (ECs/sCs) breakdown for Cdh1-Igf1r.
$data = @(
    @("ECs", "Cdh1", "Igf1r", "ECs", 1, 0.3333333333333333, 0.5650276666666666, 1.695083, 0.3735116729352592, 0.3735116729352592, 3, 1, 24.519512, 73.558536, 0.4736537296697991, 0.4736537296697991, 13.85420265316533, 124.687823878488, 0.1769151969609916, 0.1769151969609916),
    @("ECs", "Cdh1", "Igf1r", "FAPs", 1, 0.3333333333333333, 0.5650276666666666, 1.695083, 0.3735116729352592, 0.3735116729352592, 3, 1, 14.70328633333333, 44.109859, 0.284029568377475, 0.284029568377475, 8.30776356925522, 74.76987212329699, 0.1060883592477503, 0.1060883592477503),
    @("ECs", "Cdh1", "Igf1r", "sCs", 1, 0.3333333333333333, 0.5650276666666666, 1.695083, 0.3735116729352592, 0.3735116729352592, 3, 1, 12.543947, 37.631841, 0.2423167019527259, 0.2423167019527259, 7.087677104200333, 63.789093937803, 0.09050811672651723, 0.09050811672651723),
    @("sCs", "Cdh1", "Igf1r", "ECs", 3, 1, 0.9477166666666665, 2.84315, 0.6264883270647409, 0.6264883270647408, 3, 1, 24.519512, 73.558536, 0.4736537296697991, 0.4736537296697991, 23.23755018093333, 209.1379516284, 0.2967385327088075, 0.2967385327088074),
    @("sCs", "Cdh1", "Igf1r", "FAPs", 3, 1, 0.9477166666666665, 2.84315, 0.6264883270647409, 0.6264883270647408, 3, 1, 14.70328633333333, 44.109859, 0.284029568377475, 0.284029568377475, 13.93454951287222, 125.41094561585, 0.1779412091297247, 0.1779412091297247),
    @("sCs", "Cdh1", "Igf1r", "sCs", 3, 1, 0.9477166666666665, 2.84315, 0.6264883270647409, 0.6264883270647408, 3, 1, 12.543947, 37.631841, 0.2423167019527259, 0.2423167019527259, 11.88810763768333, 106.99296873915, 0.1518085852262087, 0.1518085852262087)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $rowValues.Length; $j++) {
        $col = $j + 1
        $ws.Cells.Item($row, $col).Value = $rowValues[$j]
    }
}
